$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.256752729415894
$ws.Range("B1").Value = 2.207502365112305
$ws.Range("C1").Value = 4.982283592224121
$ws.Range("D1").Value = 1.988591551780701
$ws.Range("E1").Value = 1.073990345001221
